$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "row header" cell (A2, style index 1:
# bold, bordered, centered/top-aligned) onto the new row's label cell A7,
# then set the values for the new "Overall Sentiment" row (A7:I7).
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A7").Value = "Overall Sentiment"
$ws.Range("B7").Value = "Positive"
$ws.Range("C7").Value = "Negative"
$ws.Range("D7").Value = "Negative"
$ws.Range("E7").Value = "Positive"
$ws.Range("F7").Value = "Positive"
$ws.Range("G7").Value = "Negative"
$ws.Range("H7").Value = "Positive"
$ws.Range("I7").Value = "Positive"
